# Issue #5: property boat & car done
# The "汽車" (car) sheet was generated with a bug: row 1 (meant to be the
# header row) was filled with the same values as the data row instead of
# the column labels, and the sheet was missing the common trailer columns
# (property_category .. index) that the other property sheets already
# carry. This fixes the header row and appends those columns, and adds a
# "capacity" (排氣量) column/value for the car.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 汽車 (car) sheet

# ---- Row 1: proper column headers (previously held stray data values) ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Give the new header cells (H1:N1) the same bold / centered / bordered
# look as the existing header cells (B1:G1).
$headerRange = $ws.Range("H1:N1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# ---- Row 2: data ----
# capacity (engine displacement, cc) - was mistakenly stored as text "2995"
$ws.Range("C2").Value = 2995

# Trailer columns shared with the other property sheets (land/building/...)
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# register as literal text, not an auto-converted date serial number
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2013-12-31"
$ws.Range("B2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K2").Value = "許添財"
$ws.Range("L2").Value = 639
$ws.Range("M2").Value = "tmpbb0f1"
$ws.Range("N2").Value = 32

Write-Host "sheet3 car columns fixed"
